$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> now Ricardo Andres Yara Pacheco (was Diana Maria Payares Perez)
$ws.Range("C16").Value = "20432527"
$ws.Range("D16").Value = "RICARDO ANDRES YARA PACHECO"
$ws.Range("E16").Value = "2305"
$ws.Range("F16").Value = 10000

# Row 17 -> stays Sissy Emperatriz Algarin Mendoza, only the mora value changes
$ws.Range("G17").Value = 1300000

# Row 18 -> now Diana Maria Payares Perez (was Ricardo Andres Yara Pacheco)
$ws.Range("C18").Value = "1010197164"
$ws.Range("D18").Value = "DIANA MARIA PAYARES PEREZ"
$ws.Range("E18").Value = "2404"
$ws.Range("F18").Value = 8000
